$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-10 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-11 Thursday", 2)
$d.Content.Find.Execute("66÷5=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "35÷5=7, 0", 2)
$d.Content.Find.Execute("62÷3=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "45÷4=11, 1", 2)
$d.Content.Find.Execute("57÷3=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "81÷7=11, 4", 2)
$d.Content.Find.Execute("87÷3=29, 0", $true, $false, $false, $false, $false, $true, 1, $false, "98÷9=10, 8", 2)
$d.Content.Find.Execute("23÷6=3, 5", $true, $false, $false, $false, $false, $true, 1, $false, "72÷9=8, 0", 2)
$d.Content.Find.Execute("66÷3=22, 0", $true, $false, $false, $false, $false, $true, 1, $false, "84÷7=12, 0", 2)
$d.Content.Find.Execute("50÷4=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "59÷6=9, 5", 2)
$d.Content.Find.Execute("43÷3=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "80÷8=10, 0", 2)
$d.Content.Find.Execute("43÷9=4, 7", $true, $false, $false, $false, $false, $true, 1, $false, "89÷3=29, 2", 2)
$d.Content.Find.Execute("33÷3=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "44÷3=14, 2", 2)
$d.Content.Find.Execute("15÷5=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "46÷2=23, 0", 2)
$d.Content.Find.Execute("99÷2=49, 1", $true, $false, $false, $false, $false, $true, 1, $false, "77÷5=15, 2", 2)
$d.Content.Find.Execute("78÷2=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "21÷8=2, 5", 2)
$d.Content.Find.Execute("66÷8=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "21÷2=10, 1", 2)
$d.Content.Find.Execute("88÷4=22, 0", $true, $false, $false, $false, $false, $true, 1, $false, "78÷3=26, 0", 2)
$d.Content.Find.Execute("96÷2=48, 0", $true, $false, $false, $false, $false, $true, 1, $false, "23÷2=11, 1", 2)
$d.Content.Find.Execute("28÷8=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "94÷4=23, 2", 2)
$d.Content.Find.Execute("73÷7=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "82÷4=20, 2", 2)
$d.Content.Find.Execute("81÷3=27, 0", $true, $false, $false, $false, $false, $true, 1, $false, "91÷8=11, 3", 2)
$d.Content.Find.Execute("75÷3=25, 0", $true, $false, $false, $false, $false, $true, 1, $false, "83÷3=27, 2", 2)
$d.Content.Find.Execute("33÷7=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "29÷8=3, 5", 2)
$d.Content.Find.Execute("54÷9=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "47÷2=23, 1", 2)
$d.Content.Find.Execute("69÷6=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "12÷6=2, 0", 2)
$d.Content.Find.Execute("47÷6=7, 5", $true, $false, $false, $false, $false, $true, 1, $false, "56÷9=6, 2", 2)
$d.Content.Find.Execute("52÷6=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=6, 0", 2)
